$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers - copy formatting from H1 (same header style) then set text
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Rows 2-29: I = 1, J = same as H
for ($r = 2; $r -le 29; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value()
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}

# Row 30: special values
$ws.Cells.Item(30, 9).Value = 6
$ws.Cells.Item(30, 10).Value = 7
